# Poster edits: wording tweaks, a heading restyle, and a QR-code nudge.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Figure 6 caption: "EDTs" -> "Early Decay Times" --------------------
$shp = $s.Shapes.Item(5)
$tr = $shp.TextFrame.TextRange
$oldLen = "Figure 6: EDTs for the five rooms, related to the estimated RIRs in VR environment. The dashed lines show the JND limit of 5%".Length + 2
$part = $tr.Characters(1, $oldLen)
$part.Text = "Figure 6: Early Decay Times for the five rooms, related to the estimated RIRs in VR environment. The dashed lines show the JND limit of 5%" + [char]160 + "("

# --- "Results" heading -> "Audio Evaluation Results" (+ underline) ------
$shp = $s.Shapes.Item(10)
$origHeight = $shp.Height
$tr = $shp.TextFrame.TextRange
$tr.Text = "Audio Evaluation Results"
$tr.Font.Underline = -1
$shp.Height = $origHeight

# --- Figure 7 caption: "RT60s" -> "Reverberation Times (60dB)" ----------
$shp = $s.Shapes.Item(12)
$tr = $shp.TextFrame.TextRange
$tr.Text = $tr.Text.Replace("RT60s", "Reverberation Times (60dB)")

# --- Fig 4 caption wording ------------------------------------------------
$shp = $s.Shapes.Item(16)
$tr = $shp.TextFrame.TextRange
$tr.Text = "Fig 4:Input image of Kitchen indoor room scene. "

# --- Fig 3 caption wording ------------------------------------------------
$shp = $s.Shapes.Item(18)
$tr = $shp.TextFrame.TextRange
$tr.Text = $tr.Text.Replace("only few clicks after setup", "only a few clicks provided setup")

# --- QR code picture: nudge left -----------------------------------------
$shp = $s.Shapes.Item(21)
$shp.Left = 1251.4803937007873
